$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.936.38'
$ws.Range('E2').Value = '  +2.73%  '
$ws.Range('D3').Value = '2.455.51'
$ws.Range('E3').Value = '  +2.07%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.31'
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.94'
$ws.Range('E6').Value = '  +2.89%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.541'
$ws.Range('E8').Value = '  +0.85%  '
$ws.Range('D9').Value = '2.454.39'
$ws.Range('E9').Value = '  +1.58%  '
$ws.Range('E10').Value = '  +2.37%  '
$ws.Range('E11').Value = '  +2.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.28'
$ws.Range('E12').Value = '  +1.34%  '
$ws.Range('E13').Value = '  +2.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.60'
$ws.Range('E14').Value = '  +7.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000179'
$ws.Range('E15').Value = '  +3.81%  '
$ws.Range('D16').Value = '2.900.47'
$ws.Range('E16').Value = '  +1.97%  '
$ws.Range('D17').Value = '62.845.48'
$ws.Range('E17').Value = '  +3.04%  '
$ws.Range('D18').Value = '2.460.67'
$ws.Range('E18').Value = '  +1.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.93'
$ws.Range('E19').Value = '  -1.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '330.14'
$ws.Range('E21').Value = '  +2.11%  '
$ws.Range('E22').Value = '  +1.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.15'
$ws.Range('E23').Value = '  +10.34%  '
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.43'
$ws.Range('E25').Value = '  +2.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '658.10'
$ws.Range('E26').Value = '  +7.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.19'
$ws.Range('E27').Value = '  +19.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.75'
$ws.Range('E28').Value = '  +5.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000100'
$ws.Range('E29').Value = '  +5.91%  '
$ws.Range('D30').Value = '2.576.08'
$ws.Range('E30').Value = '  +2.27%  '
$ws.Range('E31').Value = '  +2.46%  '
$ws.Range('E32').Value = '  +3.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.88'
$ws.Range('E33').Value = '  +4.15%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.139'
$ws.Range('E34').Value = '  +4.59%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.49'
$ws.Range('E35').Value = '  +1.24%  '
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('B37').Value = 'BabyDogeCoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D37').Value = '0.0₆0368'
$ws.Range('E37').Value = '  +30.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.79'
$ws.Range('E38').Value = '  +3.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.54'
$ws.Range('E39').Value = '  +4.15%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '152.84'
$ws.Range('E41').Value = '  -0.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '18.82'
$ws.Range('E42').Value = '  +2.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.75'
$ws.Range('E43').Value = '  +8.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.77'
$ws.Range('E44').Value = '  +3.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.58'
$ws.Range('E45').Value = '  +1.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '15.00'
$ws.Range('E47').Value = '  +26.99%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '146.43'
$ws.Range('E48').Value = '  +3.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.63'
$ws.Range('E49').Value = '  +2.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.61'
$ws.Range('E50').Value = '  +3.32%  '
$ws.Range('E51').Value = '  +2.34%  '
